# "se añade el flujo completo de compra"
# Adds two new rows (Cameras/Nikon D300 and Phones & PDAs/iPhone) to Sheet1's
# product catalogue, then spins off a new Sheet3 - a copy of the finished
# Sheet1 - used to flag a couple of failed purchase-flow checks ("fail") in
# column C, and leaves Sheet3 as the active/selected sheet.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# --- Sheet1: append the two new catalogue rows -----------------------------
# (write column C before column A on each row so new shared strings land in
# the same order the source workbook uses: Nikon D300, Cameras, iPhone)
$ws1.Range("C5").Value = "Nikon D300"
$ws1.Range("A5").Value = "Cameras"
$ws1.Range("C6").Value = "iPhone"
$ws1.Range("A6").Value = "Phones & PDAs"

# --- Sheet3: duplicate Sheet1 (placed after the last sheet) ----------------
$ws1.Copy([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws3 = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3.Name = "Sheet3"

# mark the two "failed" checkout rows on the copy
$ws3.Range("C3").Value = "fail"
$ws3.Range("C5").Value = "fail"

# column C on the copy re-fits to its own (now shorter) longest entry
$ws3.Columns.Item(3).ColumnWidth = 20.666666666666668

# --- view state --------------------------------------------------------------
# Sheet1: no longer the selected tab; selection becomes the whole used range
$ws1.Activate() | Out-Null
$ws1.Range("A1:C6").Select() | Out-Null

# Sheet3: becomes the active/selected tab, its own zoom level and selection
$ws3.Activate() | Out-Null
$excel.ActiveWindow.Zoom = 266
$ws3.Range("C5").Select() | Out-Null
